# "removed waste heat recovery from electricity to co2 capture.
#  Now no waste heat is considered at all."
#
# In the "connections" sheet there was a connection row describing the
# "power" chain's "simple_power" process outputting "waste heat", which was
# received as "recovered heat" feeding into the "CO2 Capture" chain's
# "simple_CO2capture" process (replacing "heat"). That whole connection row
# is removed; the remaining connection rows below it shift up to fill the
# gap (Excel does this automatically when a sheet row is deleted).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Row 22 is: power | simple_power | outflow | waste heat | recovered heat |
#            inflows | simple_CO2capture | CO2 Capture | heat
# Deleting it removes the waste-heat-recovery connection entirely and
# shifts every following row up by one.
$ws.Rows("22:22").Delete()

# Reflect the same sheet/selection focus as the edited workbook: the
# "connections" sheet becomes the active tab, with F29 selected.
$ws.Activate()
$ws.Range("F29").Select()
